# Update iServ stats for 2025-11 (row 24)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24").Value = 6408
$ws.Range("C24").Value = 1002
$ws.Range("D24").Value = 5978456
$ws.Range("E24").Value = 932.9675405742821
$ws.Range("F24").Value = 9.239686327991814
$ws.Range("G24").Value = 3.83419689119171
$ws.Range("H24").Value = 26.64685550818842
